# Restore full project from D drive for Win7 build
# Converts the monthly report sheet into a daily report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Sheet name: 月度统计_2025-12 -> 日度统计_2025-12
# ---------------------------------------------------------------------
$ws.Name = "日度统计_2025-12"

# ---------------------------------------------------------------------
# 2. Title text + remove the blank spacer row that used to sit under it
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "2025-12 日度收费统计报表"
$ws.Rows("2").Delete()

# After the delete, the layout (1-based rows) is now:
#  1  title
#  2  "统计摘要"           (was row 3)
#  3  总账单数 / 45         (was row 4)
#  4  已缴费数 / 1          (was row 5)
#  5  未缴费数 / 44         (was row 6)
#  6  总金额   / ¥107015.40 (was row 7)
#  7  已缴费金额/ ¥360.00   (was row 8)
#  8  欠费金额 / ¥106655.40 (was row 9)
#  9  缴费率   / 2.2%       (was row 10)
# 10  (blank)               (was row 11)
# 11  "收费项目明细"        (was row 12)
# 12  收费项目明细表头      (was row 13)
# 13  物业费 row             (was row 14)
# 14  生成时间 row           (was row 15)

# ---------------------------------------------------------------------
# 3. Strip the bullet-style formatting from the summary key/value rows
#    (rows 3-9) - they become plain, unstyled cells. The "统计摘要"
#    section banner shrinks from 14pt to 12pt bold.
# ---------------------------------------------------------------------
$ws.Range("A3:A9").ClearFormats()
$ws.Range("A2").Font.Size = 12

# ---------------------------------------------------------------------
# 4. Update the amended summary figures
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "¥1369.88"
$ws.Range("B8").Value = "¥105645.52"

# ---------------------------------------------------------------------
# 5. Replace the "收费项目明细" banner (row 11) with the new daily
#    breakdown table header, then insert the new data row + spacer.
# ---------------------------------------------------------------------
$ws.Range("A11").ClearFormats()
$ws.Range("A11").Value = "日期"
$ws.Range("B11").Value = "账单数"
$ws.Range("C11").Value = "日合计(¥)"
$ws.Range("D11").Value = "已缴(¥)"
$ws.Range("E11").Value = "欠费(¥)"

$ws.Rows("12:13").Insert()
$ws.Range("A12:E13").ClearFormats()

# Force the date column to remain literal text instead of being
# auto-converted to a date serial number.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-12-16"
$ws.Range("A12").ClearFormats()

$ws.Range("B12").Value = 45
$ws.Range("C12").Value = "¥107015.40"
$ws.Range("D12").Value = "¥1369.88"
$ws.Range("E12").Value = "¥105645.52"
# Row 13 stays blank (separator row).

# ---------------------------------------------------------------------
# 6. The item-detail table (now at rows 14-16) loses its header shading
#    and the summary column label "已缴费金额" becomes "已缴金额".
# ---------------------------------------------------------------------
$ws.Range("A14:F14").ClearFormats()
$ws.Range("E14").Value = "已缴金额"

# ---------------------------------------------------------------------
# 7. Update the generated-at timestamp row.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "生成时间: 2025-12-27 12:09:58"

# ---------------------------------------------------------------------
# 8. Column widths: columns A, C, D, E, F all become 18 characters wide.
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 17.15
$ws.Columns("C").ColumnWidth = 17.15
$ws.Columns("D").ColumnWidth = 17.15
$ws.Columns("E").ColumnWidth = 17.15
$ws.Columns("F").ColumnWidth = 17.15

# ---------------------------------------------------------------------
# 9. Merge range for the title shrinks from A1:D1 to A1:E1.
# ---------------------------------------------------------------------
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:E1").Merge()
